$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 618 (pushes existing row 618 and everything
# below it down by one row).
$ws.Rows.Item(618).Insert()

# Populate the newly inserted row 618 with the new data record.
$ws.Cells.Item(618, 1).Value2 = 3
$ws.Cells.Item(618, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(618, 3).Value2 = "Coquimbo"
$ws.Cells.Item(618, 4).Value2 = 45106
$ws.Cells.Item(618, 5).Value2 = 5
$ws.Cells.Item(618, 6).Value2 = 100112017
$ws.Cells.Item(618, 7).Value2 = "Apio"
$ws.Cells.Item(618, 8).Value2 = "Americana (o)"
$ws.Cells.Item(618, 9).Value2 = "Primera"
$ws.Cells.Item(618, 10).Value2 = 230
$ws.Cells.Item(618, 11).Value2 = 9000
$ws.Cells.Item(618, 12).Value2 = 9500
$ws.Cells.Item(618, 13).Value2 = 9239
$ws.Cells.Item(618, 14).Value2 = "$/docena de matas"
$ws.Cells.Item(618, 15).Value2 = "Pan de Azúcar"
$ws.Cells.Item(618, 16).Value2 = 1540
$ws.Cells.Item(618, 17).Value2 = 6
$ws.Cells.Item(618, 18).Value2 = "Hortaliza"

# Apply the same number format/style used by the other date cells in
# column D (style index 2 in the original workbook -> date/time format).
$ws.Cells.Item(618, 4).NumberFormat = $ws.Cells.Item(619, 4).NumberFormat
